# Append new temperature/humidity log rows (14-17) to the "temperatura" sheet,
# matching the data added in the upstream commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @("2024-11-27 18:09:29", 29, "Alerta Umidade nao encontrado"),
    @("2024-11-27 18:23:12", 29, "Alerta Umidade nao encontrado"),
    @("2024-11-27 18:45:37", 31, "Alerta Umidade nao encontrado"),
    @("2024-11-27 18:48:50", 30, "Alerta Umidade nao encontrado")
)

$startRow = 14
for ($i = 0; $i -lt $newRows.Length; $i++) {
    $rowNum = $startRow + $i
    $rowData = $newRows[$i]

    $ws.Cells.Item($rowNum, 1).Value = $rowData[0]
    $ws.Cells.Item($rowNum, 2).Value = $rowData[1]
    $ws.Cells.Item($rowNum, 3).Value = $rowData[2]
}
